# Add data for 2022-02-24
# - Rename sheet/title from "Through 2022-02-15" to "Through 2022-02-16"
# - Update the "I" column header text accordingly
# - Update the two affected data cells (I3 and I14)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (flows through to <sheet name="..."> in workbook.xml)
$ws.Name = "Through 2022-02-16"

# Update the header label in column I (shared string used by I1)
$ws.Range("I1").Value = "2022 (through 02-16)"

# Update the updated counts
$ws.Range("I3").Value = 77
$ws.Range("I14").Value = 238
